$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student record appended as row 50 (STT 49).
# Set the numeric cells first (while the cell still carries the default
# style) so the later format-copy doesn't get confused about their type.
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 19525742

# Copy the formatting (styles/number formats) of the last existing data
# row (49) down onto the new row 50, matching how the sheet's existing
# rows are all formatted identically.
$ws.Range("A49:H49").Copy()
$ws.Range("A50:H50").PasteSpecial(-4122)

# Now fill in the text cells. Column F/G hold date-like / numeric-looking
# text (must stay text so the leading zero on the phone number survives),
# so enter them with a leading apostrophe just like Excel's own
# text-quote-prefix entry.
$ws.Range("C50").Value = "Trần Nguyễn Kha"
$ws.Range("D50").Value = "Hoàng"
$ws.Range("E50").Value = "Nam"
$ws.Range("F50").Value = "'" + "30/04/2002"
$ws.Range("G50").Value = "'" + "0972211750"
$ws.Range("H50").Value = "DHKTPM15C"

# Match the author's final on-screen view/selection.
$ws.Range("I50").Select()
